# Penalty Reward System update (unfinished, per commit message):
# the whole 16-week forecast window rolls forward by one week, the
# MyForecast numbers are refreshed, and the Summary sheet's derived
# figures are recomputed to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet: Week_Start_Date (B) & MyForecast (D) ---
# Week_Start_Date cells hold plain text like "2025-01-12" in the source
# file; format as Text first so COM doesn't silently turn the string into
# a date serial number.

$forecastRows = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 39 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 40 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 41 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 43 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 23 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 24 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 25 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 37 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 32 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 31 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 23 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 24 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 34 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 31 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 30 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 31 }
)

foreach ($entry in $forecastRows) {
    $dateCell = $wsForecast.Cells.Item($entry.Row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $entry.Date

    $wsForecast.Cells.Item($entry.Row, 4).Value = $entry.Forecast
}

# --- "Summary" sheet: recomputed metrics ---
# Every value in column B of this sheet is stored as text, including the
# ones that look numeric, so force Text formatting before assigning.

$summaryRows = @(
    @{ Row = 2;  Value = "2023-02-19 to 2025-01-05" },
    @{ Row = 8;  Value = "1573 units" },
    @{ Row = 9;  Value = "508" },
    @{ Row = 10; Value = "272" },
    @{ Row = 11; Value = "163" },
    @{ Row = 12; Value = "43" },
    @{ Row = 13; Value = "2025-02-02" },
    @{ Row = 14; Value = "23" },
    @{ Row = 15; Value = "2025-02-09" }
)

foreach ($entry in $summaryRows) {
    $cell = $wsSummary.Cells.Item($entry.Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $entry.Value
}
